$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range('A2').Value = 'Lancement du jeu'
$ws.Range('B2').Value = 'En tant que joueur 
Je veux que le jeux se lance sur la page du menu
Pour pouvoir choisir ce que je veux faire'
$ws.Range('C2').Value = 'Quand le jeu n''est pas lancé sur mon pc
Au lancement de l''éxecutable
Le jeu se lance sur le menu principale '
$ws.Range('D2').Value = 'Quand le jeu n''est pas lancé sur mon pc
Au lancement de l''éxecutable
Le curseur se trouve sur l''entrée de tout en haut'

$ws.Range('A3').Value = 'Menu principale'
$ws.Range('B3').Value = 'En tant que joueur,
Je veux utiliser les flèches haut/bas
Pour sélectionner une entrée dans un menu'
$ws.Range('C3').Value = 'Quand je suis sur le menu principale
A la pression des flèches haut et bas
Le curseur de séléction change d''entrée d''une ligne dans la direction choisie'
$ws.Range('D3').Value = 'Quand le curseur sélectionne l''entrée de tout en haut
A la pression de la flèche haut
Il ne se passe rien'
$ws.Range('E3').Value = 'Quand le curseur sélectionne l''entrée de tout en bas
A la pression de la flèche bas
Il ne se passe rien'
$ws.Range('F3').Value = 'Quand le curseur pointe sur une entrée
A la pression de la touche entrée
L''entrée sélectionnée s''affiche'

$ws.Range('A5').Value = 'Menu option'
$ws.Range('B5').Value = 'En tant que joueur,
Je veux accéder au menu option
afin de pouvoir modifier la difficulté et le son du jeu'
$ws.Range('C5').Value = 'dans le menu option avec le curseur sur l''option de difficulté sélectionnée,
a la pression de la flèche droite
la difficulté change (facile -> moyen -> difficile -> godmod -> facile)'
$ws.Range('D5').Value = 'dans le menu option avec le curseur sur l''option de son sélectionnée,
a la pression de la flèche de droite
le son change (activé -> désactivé -> activé)'
$ws.Range('E5').Value = 'En tant que joueur,
Je veux utiliser les flèches haut/bas
Pour sélectionner une entrée dans le menu option'
$ws.Range('F5').Value = 'Quand je suis sur le menu option
A la pression des flèches haut et bas
Le curseur de séléction change d''entrée d''une ligne dans la direction choisie'
$ws.Range('G5').Value = 'Quand le curseur sélectionne l''entrée de tout en haut
A la pression de la flèche haut
Il ne se passe rien'
$ws.Range('H5').Value = 'Quand le curseur sélectionne l''entrée de tout en bas
A la pression de la flèche bas
Il ne se passe rien'

$ws.Range('C6').Value = 'dans le menu option avec le curseur sur l''option de difficulté sélectionnée,
a la pression de la flèche gauche
la difficulté change (facile ->  godmod -> difficile -> moyen -> facile)'
$ws.Range('D6').Value = 'dans le menu option avec le curseur sur l''option de son sélectionnée,
a la pression de la flèche de gauche
le son change (activé -> désactivé -> activé)'

$ws.Range('A7').Value = 'Pseudo avant de jouer'
$ws.Range('B7').Value = 'En tant que joueur,
Je veux rentrer mon pseudo avant de jouer
Afin que mon highscore s''affiche dans le menu "HIGHSCORES"'
$ws.Range('C7').Value = 'quand l''entrée "PLAY" sur le menu principale est choisi
une page d''entrée de pseudo s''affichage
Pour pouvoir rentrer mon pseudo
'
$ws.Range('D7').Value = 'En tant que joueur sur la page d''entrée de pseudo
En pressant entrée le pseudo s''enregistre
Pour savoir qui vas jouer'

$ws.Range('A8').Value = 'Lancement de la partie'
$ws.Range('B8').Value = 'En tant que joueur
Je veux que tout les éléments du jeu s''affiche correctement
Pour pouvoir jouer dans de bonnes condition'
$ws.Range('C8').Value = 'En tant que joueur au lancement de la partie
le compteur de vie s''affiche avec toute les vies
Pour savoir combien de vie il me reste
'
$ws.Range('D8').Value = 'En tant que joueur au lancement de la partie
Les murs s''affichent entierement construit
Pour pouvoir me protéger des aliens'
$ws.Range('E8').Value = 'En tant que joueur au lancement de la partie
Le vaisseau du joueur s''affiche 
Pour pouvoir tirer sur les aliens'
$ws.Range('F8').Value = 'En tant que joueur au lancement de la partie
Le compteur de score s''affiche
pour savoir mon score en temps réel'
$ws.Range('G8').Value = 'En tant que joueur au lancement de la partie
les aliens s''affichent
Pour pouvoir tirer sur les aliens et qu''ils nous tirent dessus'

$ws.Range('A9').Value = 'Tir des aliens'
$ws.Range('B9').Value = 'En tant que joueur
Je veux que les aliens me tirent dessus
Pour avoir un adversaire digne de ce nom'
$ws.Range('C9').Value = 'En tant que joueur quand la partie est lancée
Les aliens tirent vers le bas à interval irrégulier
Pour que les aliens me tirent dessus'

$ws.Range('A10').Value = 'Déplacement des aliens'
$ws.Range('B10').Value = 'En tant que joueur 
Je veux que les aliens se déplacent latéralement et verticalement
Pour ajouter un atout a mon adversaire'
$ws.Range('C10').Value = 'Quand la partie est en cours
Les aliens se déplacent vers la droite a un rythme régulier
Pour faires bouger les aliens
'
$ws.Range('D10').Value = 'Quand la partie est en cours
Les aliens se déplacent vers la gauche a un rythme régulier
Pour faires bouger les aliens
'
$ws.Range('E10').Value = 'Quand les aliens arrivent au bord de l''écran
Les aliens descende d''une ligne et reparte dans l''autre sens
Pour que les aliens se rapproche petit à petit du joueur'
$ws.Range('F10').Value = 'Quand les aliens se déplacent horizontalement 
A chaque case parcourue le model visuel de l''alien change
Afin de créer une petite animation de déplacement'

$ws.Range('A11').Value = 'Compteur de vie'
$ws.Range('B11').Value = 'En tant que joueur
je veux avoir un compteur de vie ( 3 vie )
pour savoir combien de fois j''ai encore droit a l''erreur'
$ws.Range('C11').Value = 'Quand un tir alien touche le vaisseau du joueur
le compteur fait disparaitre une vie
Pour que le joueur perde une vie'
$ws.Range('D11').Value = '
'

$ws.Range('A12').Value = 'Game over '
$ws.Range('B12').Value = 'En tant que joueur 
Je veux savoir quand j''ai perdu la partie'
$ws.Range('C12').Value = 'Quand le compteur de vie atteint zero
Le jeu s''arrète et un écran de Game Over s''affiche
Pour savoir quand j''ai perdu'
$ws.Range('D12').Value = 'Quand les aliens atteinent la ligne où commence les murs de protection
Le jeu s''arrète et un écran de Game Over s''affiche
Pour savoir quand j''ai perdu
'

$ws.Range('A13').Value = 'vaisseau du joueur touché'
$ws.Range('B13').Value = 'En tant que joueur
Je veux savoir quand un tir ennemi me touche
Pour savoir quand je pert une vie'
$ws.Range('C13').Value = 'Quand le vaisseau du joueur est touché
il s''immobilise 2s et devient intouchable
Pour pénalisé le joueur sans lui faire perdre d''autre vie'
$ws.Range('D13').Value = 'Quand le vaisseau du joueur est touché
Le vaisseau change de forme durant 2s
Pour donner un signal visuel au joueur
'

$ws.Range('A14').Value = 'Scores en temps réel'
$ws.Range('B14').Value = 'En tant que joueur
je veux que chaque alien détruit me rapporte des points
Afin de battre mon highscore'
$ws.Range('C14').Value = 'Quand un alien est détruit 
il rapporte des points
pour pouvoir augmenté mon score'
$ws.Range('D14').Value = 'Quand le vaisseau du joueur est touché
Le joueur pert des points
Pour pénalisé les joueur qui se font beaucoup touché
'

$ws.Range('A15').Value = 'Aliens détruit'
$ws.Range('B15').Value = 'En tant que joueur
Je veux pouvoir détruire les aliens
Pour pouvoir gagner des points'
$ws.Range('C15').Value = 'Quand un alien est touché par un tir du vaisseau du joueur
l''alien disparait dans un petite explosion ( animation )
Pour donner un signal visuel au joueur'

$ws.Range('A16').Value = 'Tout les aliens détruit'
$ws.Range('B16').Value = 'En tant que joueur
Je veux que le jeu réagisse a la mort d''une vague entière d''aliens
Pour pouvoir continuer à jouer et être récompensé'
$ws.Range('C16').Value = 'Quand tout les aliens sont détruit
le joueur a 3s de répit avant qu''une nouvelle vague d''aliens apparaisse
Pour que la partie se continue'
$ws.Range('D16').Value = 'Quand tout les aliens sont détuit
Le joueur recois une vie supplémentaire ( dans la limite maximal du compteur de vie, 3 vie )
Pour récompensé le joueur d''avoir éliminé un vague entière d''énnemis'

$ws.Range('A17').Value = 'Déplacement du vaisseau du joueur'
$ws.Range('B17').Value = 'En tant que joueur 
Je veux pouvoir déplacer mon vaisseau de manière horizontal
Pour pouvoir esquiver les tirs ennemis et me positionné pour tirer  sur les ennemis'
$ws.Range('C17').Value = 'Quand la flèche de droite est préssée en jeu
Le vaisseau se déplace sur la droite
Pour pouvoir me déplacer
'
$ws.Range('D17').Value = 'Quand la flèche de gauche est préssée en jeu
Le vaisseau se déplace sur la gauche
Pour pouvoir me déplacer
'
$ws.Range('E17').Value = 'Quand une flèche directionel est appuyée longuement
Le vaisseau se déplace de manière continue dans la direction séléctionée
Pour facilité les déplacement de joueur'
$ws.Range('F17').Value = 'Quand le vaisseau arrive au bord de l''écran
Le vaisseau ne peut plus se déplacer vers le côté ou il est déjà collé au bord
Afin que le vaisseau ne disparaisse pas de l''écran'

$ws.Range('A18').Value = 'Tir du vaisseau du joueur'
$ws.Range('B18').Value = 'En tant que joueur 
Je veux pouvoir tirer vers le haut avec mon vaisseau
Pour pouvoir détruire les aliens'
$ws.Range('C18').Value = 'Quand la flèche haut est préssée en jeu
Le vaisseau tir un laser vers le haut
Pour pouvoir tirer sur les aliens
'
$ws.Range('D18').Value = 'Quand un laser est tirer par le vaisseau du joueur
Le vaisseau ne peut pas tirer d''autre laser pendant 1s
Pour limiter le nombre de tir du joueur
'

$ws.Range('A19').Value = 'Tir sur les murs'
$ws.Range('B19').Value = 'En tant que joueur 
Je veux que les murs soit déstruictible
Pour offrir plus de possibilité de gameplay'
$ws.Range('C19').Value = 'Quand un tir touche un mur
la partie du mur change de couleur ( blanc -> rouge )
Pour indiquer au joueur que le mur vas bientôt casser'
$ws.Range('D19').Value = 'Quand un tir touche un mur endommagé
La partie du mur disparait
Afin de laisser passer les tirs
'
